$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basketball")

$ws.Range("C3").Value = 0.1
$ws.Range("C4").Value = 0.1
$ws.Range("C5").Value = 0.1
